# Update tracklist from "The Weeknd - After Hours" (weeknd6) to
# "R+R=Now" (rplusrequalnow1): new track data, renamed defined names,
# and Sheet2 selection follow-up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rewrite the track data on Sheet1 and Sheet3 (they mirror each other)
# ---------------------------------------------------------------------
$sheetNames = @("Sheet1", "Sheet3")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2 - Change of Tone
    $ws.Range("B2").Value = "Change of Tone"
    $ws.Range("C2").Value = "Christian Scott aTunde Adjuah / Robert Glasper / Derrick Hodge / Terrace Martin / Taylor McFerrin / Rose McKinney / Justin Tyson"
    $ws.Range("D2").Value = "R+R=Now"
    $ws.Range("E2").Value = 0.31388888888888888

    # Row 3 - Awake to You
    $ws.Range("B3").Value = "Awake to You"
    $ws.Range("C3").Value = "Christian Scott aTunde Adjuah / Robert Glasper / Derrick Hodge / Terrace Martin / Taylor McFerrin / Rose McKinney / Justin Tyson"
    $ws.Range("D3").Value = "R+R=Now"
    $ws.Range("E3").Value = 0.3430555555555555

    # Row 4 - By Design
    $ws.Range("B4").Value = "By Design"
    $ws.Range("C4").Value = "Christian Scott aTunde Adjuah / Robert Glasper / Derrick Hodge / Jahi Lake / Terrace Martin / Taylor McFerrin / Rose McKinney / Justin Tyson"
    $ws.Range("D4").Value = "R+R=Now"
    $ws.Range("E4").Value = 0.22708333333333333

    # Row 5 - Resting Warrior
    $ws.Range("B5").Value = "Resting Warrior"
    $ws.Range("C5").Value = "Christian Scott aTunde Adjuah / Robert Glasper / Derrick Hodge / Terrace Martin / Rose McKinney / Justin Tyson"
    $ws.Range("D5").Value = "R+R=Now"
    $ws.Range("E5").Value = 0.41041666666666665

    # Row 6 - Needed You Still (B6 keeps its s="7" style, untouched by value write)
    $ws.Range("B6").Value = "Needed You Still"
    $ws.Range("C6").Value = "Christian Scott aTunde Adjuah / Robert Glasper / Omari Hardwick / Derrick Hodge / Terrace Martin / Taylor McFerrin / Rose McKinney / Michael Neil / Justin Tyson"
    $ws.Range("D6").Value = "R+R=Now feat. Omari Hardwick"
    $ws.Range("E6").Value = 0.25625000000000003

    # Row 7 - Escape from LA / same composer text as before, new performer+time
    $ws.Range("D7").Value = "R+R=Now"
    $ws.Range("E7").Value = 0.19999999999999998

    # Row 8 - Heartless / same composer text as before, new performer+time
    $ws.Range("D8").Value = "R+R=Now feat. Terry Crews"
    $ws.Range("E8").Value = 0.32013888888888892

    # Row 9 - Reflect (Reprise)
    $ws.Range("B9").Value = "Reflect (Reprise)"
    $ws.Range("C9").Value = "Robert Glasper / Derrick Hodge / Terrace Martin / Taylor McFerrin / Rose McKinney / Kyle Myricks / Justin Tyson"
    $ws.Range("D9").Value = "R+R=Now feat. Stally"
    $ws.Range("E9").Value = 0.27916666666666667

    # Row 10 - Her=Now
    $ws.Range("B10").Value = "Her=Now"
    $ws.Range("C10").Value = "Derrick Hodge / Taylor McFerrin / Amanda Seales"
    $ws.Range("D10").Value = "R+R=Now feat. Amanda Seales"
    $ws.Range("E10").Value = 0.28055555555555556

    # Row 11 - Respond
    $ws.Range("B11").Value = "Respond"
    $ws.Range("C11").Value = "Christian Scott aTunde Adjuah / Derrick Hodge / Taylor McFerrin / Justin Tyson"
    $ws.Range("D11").Value = "R+R=Now"
    $ws.Range("E11").Value = 0.22291666666666665

    # Row 12 - Been on My Mind
    $ws.Range("B12").Value = "Been on My Mind"
    $ws.Range("C12").Value = "Christian Scott aTunde Adjuah / Robert Glasper / Derrick Hodge / Jahi Lake / Taylor McFerrin / Amber Navran / Dante Smith"
    $ws.Range("D12").Value = "R+R=Now feat. Amber Navran"
    $ws.Range("E12").Value = 0.20069444444444443

    # Rows 13-15 previously held tracks 12-14; the new tracklist only has
    # 11 tracks, so these rows are cleared out (B13/C13 and E13:E15 keep
    # their existing cell styles, just with no content).
    $ws.Range("A13:A15").ClearContents()
    $ws.Range("B13:B15").ClearContents()
    $ws.Range("C13:C15").ClearContents()
    $ws.Range("D13:D15").ClearContents()
}

# ---------------------------------------------------------------------
# 2. Rename the defined names (weeknd6 -> rplusrequalnow1) and shrink
#    their ranges from 15 to 12 rows to match the new track count.
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.RefersTo -like "*Sheet1*") {
        $n.Name = "rplusrequalnow1"
        $n.RefersTo = "=Sheet1!`$A`$1:`$E`$12"
    } elseif ($n.RefersTo -like "*Sheet3*") {
        $n.Name = "rplusrequalnow1"
        $n.RefersTo = "=Sheet3!`$A`$1:`$E`$12"
    }
}

# ---------------------------------------------------------------------
# 3. Sheet2's formula table now only needs to show through row 15
#    (3 header/rule rows + 11 tracks -> ends K15 instead of K18).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A3:K15").Select()
